$wb = $excel.ActiveWorkbook

# Sheet ALC, row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 520.2692
$ws.Range("I92").Value = 509.35
$ws.Range("J92").Value = 556.6667
$ws.Range("K92").Value = 509.35
$ws.Range("L92").Value = 556.6667
$ws.Range("M92").Value = 738.65
$ws.Range("N92").Value = -3052.6667

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1620
$ws.Range("I100").Value = 800
$ws.Range("J100").Value = 2166.6667
$ws.Range("K100").Value = 800
$ws.Range("L100").Value = 2166.6667
$ws.Range("M100").Value = -259
$ws.Range("N100").Value = -3248.6667

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1019.3871
$ws.Range("J112").Value = 1036.7
$ws.Range("L112").Value = 3110.1
$ws.Range("N112").Value = -5326.1

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5325268
$ws.Range("I132").Value = 5688104.5
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 17064313.5
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -17061783.5
$ws.Range("N132").Value = -16060.0001

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4352.1113
$ws.Range("I141").Value = 4481.2856
$ws.Range("K141").Value = 13443.8568
$ws.Range("M141").Value = -8263.856800000001

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 63808.06
$ws.Range("I45").Value = 113258.22
$ws.Range("J45").Value = 8176.625
$ws.Range("K45").Value = 113258.22
$ws.Range("L45").Value = 8176.625
$ws.Range("M45").Value = -112881.22
$ws.Range("N45").Value = -8930.625

# Sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2364.2222
$ws.Range("I63").Value = 1846.3334
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 1846.3334
$ws.Range("L63").Value = 3400
$ws.Range("M63").Value = -1160.3334
$ws.Range("N63").Value = -4772

# Sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2364.2222
$ws.Range("I66").Value = 1846.3334
$ws.Range("J66").Value = 3400
$ws.Range("K66").Value = 9231.666999999999
$ws.Range("L66").Value = 17000
$ws.Range("M66").Value = -5799.666999999999
$ws.Range("N66").Value = -23864

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 128382.375
$ws.Range("I102").Value = 254269.75
$ws.Range("J102").Value = 2495
$ws.Range("K102").Value = 254269.75
$ws.Range("L102").Value = 2495
$ws.Range("M102").Value = -252647.75
$ws.Range("N102").Value = -5739

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20950.258
$ws.Range("I132").Value = 23902.654
$ws.Range("K132").Value = 71707.962
$ws.Range("M132").Value = -69177.962

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 126984.94
$ws.Range("I105").Value = 113003.336
$ws.Range("J105").Value = 144961.28
$ws.Range("K105").Value = 113003.336
$ws.Range("L105").Value = 144961.28
$ws.Range("M105").Value = -111256.336
$ws.Range("N105").Value = -148455.28

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 83371640
$ws.Range("I107").Value = 111162060
$ws.Range("J107").Value = 371
$ws.Range("K107").Value = 111162060
$ws.Range("L107").Value = 371
$ws.Range("M107").Value = -111160140
$ws.Range("N107").Value = -4211

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12521.282
$ws.Range("I134").Value = 14063.846
$ws.Range("J134").Value = 3927
$ws.Range("K134").Value = 42191.538
$ws.Range("L134").Value = 11781
$ws.Range("M134").Value = -39656.538
$ws.Range("N134").Value = -16851

# Sheet CRP, row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 6666
$ws.Range("I10").Value = 6666
$ws.Range("K10").Value = 6666
$ws.Range("M10").Value = -6527

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3348.1667
$ws.Range("I86").Value = 2814.6365
$ws.Range("J86").Value = 3799.6155
$ws.Range("K86").Value = 2814.6365
$ws.Range("L86").Value = 3799.6155
$ws.Range("M86").Value = -1691.6365
$ws.Range("N86").Value = -6045.6155

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3348.1667
$ws.Range("I89").Value = 2814.6365
$ws.Range("J89").Value = 3799.6155
$ws.Range("K89").Value = 14073.1825
$ws.Range("L89").Value = 18998.0775
$ws.Range("M89").Value = -8457.182500000001
$ws.Range("N89").Value = -30230.0775

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1410.6666
$ws.Range("I105").Value = 1296.8
$ws.Range("K105").Value = 1296.8
$ws.Range("M105").Value = 450.2

# Sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1104
$ws.Range("J34").Value = 1271.4286
$ws.Range("L34").Value = 3814.2858
$ws.Range("N34").Value = -3982.2858

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58886684
$ws.Range("I80").Value = 125130744
$ws.Range("J80").Value = 3074
$ws.Range("K80").Value = 125130744
$ws.Range("L80").Value = 3074
$ws.Range("M80").Value = -125129746
$ws.Range("N80").Value = -5070

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 58886684
$ws.Range("I83").Value = 125130744
$ws.Range("J83").Value = 3074
$ws.Range("K83").Value = 625653720
$ws.Range("L83").Value = 15370
$ws.Range("M83").Value = -625648728
$ws.Range("N83").Value = -25354

# Sheet GSM, row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 46495
$ws.Range("J137").Value = 46495
$ws.Range("L137").Value = 46495
$ws.Range("N137").Value = -56695

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2754.12
$ws.Range("I68").Value = 1886.6666
$ws.Range("K68").Value = 1886.6666
$ws.Range("M68").Value = -1137.6666

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2754.12
$ws.Range("I71").Value = 1886.6666
$ws.Range("K71").Value = 9433.333000000001
$ws.Range("M71").Value = -5689.333000000001

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2107.5557
$ws.Range("I100").Value = 1820
$ws.Range("J100").Value = 2467
$ws.Range("K100").Value = 1820
$ws.Range("L100").Value = 2467
$ws.Range("M100").Value = -1279
$ws.Range("N100").Value = -3549

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2766.1428
$ws.Range("I122").Value = 2672.4
$ws.Range("J122").Value = 3000.5
$ws.Range("K122").Value = 8017.200000000001
$ws.Range("L122").Value = 9001.5
$ws.Range("M122").Value = -5567.200000000001
$ws.Range("N122").Value = -13901.5

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6995324.5
$ws.Range("I62").Value = 19232418
$ws.Range("K62").Value = 19232418
$ws.Range("M62").Value = -19231794

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6995324.5
$ws.Range("I65").Value = 19232418
$ws.Range("K65").Value = 96162090
$ws.Range("M65").Value = -96158970

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3180.2258
$ws.Range("I132").Value = 3150.36
$ws.Range("K132").Value = 9451.08
$ws.Range("M132").Value = -6921.08

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1948.6842
$ws.Range("I136").Value = 764.8
$ws.Range("J136").Value = 3264.111
$ws.Range("K136").Value = 2294.4
$ws.Range("L136").Value = 9792.332999999999
$ws.Range("M136").Value = 255.6000000000004
$ws.Range("N136").Value = -14892.333
